$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 114.8270160096505
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 130.6186648591143

# Row 3
$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 10.29869402782916
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 45.1914812416419
